$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.3327402135231317
$ws1.Range("C2").Value = 0.06516290726817042
$ws1.Range("D2").Value = 0.9285714285714286
$ws1.Range("E2").Value = 0.1217798594847775
$ws1.Range("F2").Value = 0.2544031311154599
$ws1.Range("G2").Value = 0.6151046405823476
$ws1.Range("H2").Value = 0.8103263777421081
$ws1.Range("I2").Value = 26
$ws1.Range("J2").Value = 373
$ws1.Range("K2").Value = 161
$ws1.Range("L2").Value = 2

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.9877300613496932
$ws2.Range("C2").Value = 0.301498127340824
$ws2.Range("D2").Value = 0.4619799139167862

$ws2.Range("B3").Value = 0.06516290726817042
$ws2.Range("C3").Value = 0.9285714285714286
$ws2.Range("D3").Value = 0.1217798594847775

$ws2.Range("B4").Value = 0.3327402135231317
$ws2.Range("C4").Value = 0.3327402135231317
$ws2.Range("D4").Value = 0.3327402135231317
$ws2.Range("E4").Value = 0.3327402135231317

$ws2.Range("B5").Value = 0.5264464843089318
$ws2.Range("C5").Value = 0.6150347779561263
$ws2.Range("D5").Value = 0.2918798867007819

$ws2.Range("B6").Value = 0.9417658615022153
$ws2.Range("C6").Value = 0.3327402135231317
$ws2.Range("D6").Value = 0.4450304450127003

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 161
$ws3.Range("C2").Value = 373

$ws3.Range("B3").Value = 2
$ws3.Range("C3").Value = 26
